$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ProgramModule")
$ws.Activate()

$ws.Range("A2").Value = "Team23JavaS"
$ws.Range("B2").Value = "AutomationTeam23cod"

$ws.Range("B6").Select()
